$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Defined names: repoint TbCidades_Cidade at the local table (it used to
#    point at an external workbook reference), add the two new named ranges,
#    and break the now-unused external link.
# ---------------------------------------------------------------------------
$nmCidade = $wb.Names.Item("TbCidades_Cidade")
$nmCidade.RefersTo = "=TbCidades[Cidade]"

$links = $wb.LinkSources()
foreach ($l in $links) {
    $wb.BreakLink($l, 1)
}

$wb.Names.Add("TbCidades_Estado", "=TbCidades[Estado]")
$wb.Names.Add("TbEstados_Estado", "=TbEstados[Estado]")

# ---------------------------------------------------------------------------
# 2. "Cidades" sheet: new helper formulas + data validation lists that drive
#    the dependent dropdowns.
# ---------------------------------------------------------------------------
$wsCidades = $wb.Worksheets.Item("Cidades")

$wsCidades.Range("D6").Formula = "=MATCH(D3,TbCidades_Estado,0)"
$wsCidades.Range("D9").Formula = "=COUNTIF(TbCidades_Estado,D3)"

$wsCidades.Range("F2").Validation.Add(3, 1, 1, "=OFFSET(TbCidades_Cidade,4,0,2)")
$wsCidades.Range("F4").Validation.Add(3, 1, 1, "=OFFSET(TbCidades_Cidade,D6-1,0,3)")
$wsCidades.Range("F6").Validation.Add(3, 1, 1, "=OFFSET(TbCidades_Cidade,D6-1,0,D9)")

# ---------------------------------------------------------------------------
# 3. "Repasses" sheet: two new rows of data, grow the table to match, and add
#    the validation dropdowns for Estado / Cidade.
# ---------------------------------------------------------------------------
$wsRepasses = $wb.Worksheets.Item("Repasses")

$wsRepasses.Range("A5").Value = "Paraná"
$wsRepasses.Range("B5").Value = "Curitiba"
$wsRepasses.Range("C5").Value = 30000
$wsRepasses.Range("A6").Value = "Pernambuco"
$wsRepasses.Range("B6").Value = "Recife"
$wsRepasses.Range("C6").Value = 50000

$tblRepasses = $wsRepasses.ListObjects.Item("TbRepasses")
$tblRepasses.Resize($wsRepasses.Range("A1:C6"))

$wsRepasses.Range("A2:A6").Validation.Add(3, 1, 1, "=TbEstados_Estado")
$wsRepasses.Range("B2:B6").Validation.Add(3, 1, 1, "=OFFSET(TbCidades_Cidade, MATCH(A2,TbCidades_Estado,0)-1,0, COUNTIF(TbCidades_Estado,A2))")

# ---------------------------------------------------------------------------
# 4. "Testes" sheet: change the dropdown pick and wire up two new formulas.
# ---------------------------------------------------------------------------
$wsTestes = $wb.Worksheets.Item("Testes")

$wsTestes.Range("F4").Value = "Médio"
$wsTestes.Range("D4").Formula = "=SUM(OFFSET(A5,0,1,3))"
$wsTestes.Range("F7").Formula = "=MATCH(F4,A2:A11,)"

# ---------------------------------------------------------------------------
# 5. Restore the view/selection state recorded in the saved workbook: visit
#    each sheet in order, leaving "Repasses" as the active tab.
# ---------------------------------------------------------------------------
$wsEstados = $wb.Worksheets.Item("Estados")
$wsEstados.Activate()
$wsEstados.Range("A2:A6").Select()

$wsCidades.Activate()
$wsCidades.Range("F6").Select()

$wsTestes.Activate()
$wsTestes.Range("F8").Select()

$wsRepasses.Activate()
$wsRepasses.Range("C7").Select()
